# Update the "Binary Search 2" sheet: row 7/8 "Binary S2" labels bump from 4 -> 6
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Binary Search 2")

$ws.Range("C8").Value = "Binary S2 6, Notes 1, Notes 2, Readme.md 3rd notes"
$ws.Range("C7").Value = "Binary S2 6"

$ws.Activate()
$ws.Range("C8").Select()
